$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the number of hours (D2) from 15 to 14; dependent formulas
# (G3, B12, B13) recalc automatically.
$ws.Range("D2").Value = 14

# Update the active selection to H13.
$ws.Range("H13").Select()
